# Update the "handback status" timestamps to reflect a regenerated report.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the bab163ec row (row 2)
$wsOverview.Range("G2").Value = "2016-08-30 17:17:37"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the bab163ec row (row 2)
$wsZhCn.Range("H2").Value = "2016-08-30 17:17:32"
$wsZhCn.Range("K2").Value = "2016-08-30 17:17:50"

# de-de sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the bab163ec row (row 2)
$wsDeDe.Range("H2").Value = "2016-08-30 17:17:37"
$wsDeDe.Range("K2").Value = "2016-08-30 17:17:57"
